# Edit script: prepend 11 new daily price rows (2022-09-13 .. 2022-09-27) to the
# "Output" sheet, shifting existing history down and dropping the oldest 11 rows
# so the sheet stays at 101 total rows (header + 100 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push the existing data rows (2:101) down by 11 rows to make room for the
#    new rows at the top. This is a native row-insert, so existing cell
#    styles/values/types are preserved exactly (no re-typing/coercion risk).
$ws.Rows("2:12").Insert()

# 2) The insert pushed the old trailing rows (previously 91:101, the oldest
#    11 days) down to 102:112. Drop them so the sheet stays at 101 rows total.
$ws.Rows("102:112").Delete()

# 3) The newly inserted rows (2:12) are blank. Force them to behave like the
#    existing text-typed data cells: mark as Text *before* assigning values so
#    numeric-looking strings (dates, prices, volumes) are kept as literal text
#    instead of being auto-converted to numbers/dates.
$ws.Range("A2:F12").NumberFormat = "@"


$ws.Cells.Item(2, 1).Value = "2022-09-27"
$ws.Cells.Item(2, 2).Value = "123.8700"
$ws.Cells.Item(2, 3).Value = "124.0700"
$ws.Cells.Item(2, 4).Value = "120.6000"
$ws.Cells.Item(2, 5).Value = "121.6100"
$ws.Cells.Item(2, 6).Value = "1040449"

$ws.Cells.Item(3, 1).Value = "2022-09-26"
$ws.Cells.Item(3, 2).Value = "123.1800"
$ws.Cells.Item(3, 3).Value = "124.2500"
$ws.Cells.Item(3, 4).Value = "121.6300"
$ws.Cells.Item(3, 5).Value = "122.3100"
$ws.Cells.Item(3, 6).Value = "1733060"

$ws.Cells.Item(4, 1).Value = "2022-09-23"
$ws.Cells.Item(4, 2).Value = "123.7300"
$ws.Cells.Item(4, 3).Value = "124.4300"
$ws.Cells.Item(4, 4).Value = "121.4400"
$ws.Cells.Item(4, 5).Value = "123.4800"
$ws.Cells.Item(4, 6).Value = "1481866"

$ws.Cells.Item(5, 1).Value = "2022-09-22"
$ws.Cells.Item(5, 2).Value = "125.5500"
$ws.Cells.Item(5, 3).Value = "126.0400"
$ws.Cells.Item(5, 4).Value = "123.7628"
$ws.Cells.Item(5, 5).Value = "124.3900"
$ws.Cells.Item(5, 6).Value = "1374943"

$ws.Cells.Item(6, 1).Value = "2022-09-21"
$ws.Cells.Item(6, 2).Value = "128.7300"
$ws.Cells.Item(6, 3).Value = "130.5750"
$ws.Cells.Item(6, 4).Value = "126.3250"
$ws.Cells.Item(6, 5).Value = "126.4800"
$ws.Cells.Item(6, 6).Value = "1105658"

$ws.Cells.Item(7, 1).Value = "2022-09-20"
$ws.Cells.Item(7, 2).Value = "129.1000"
$ws.Cells.Item(7, 3).Value = "129.6050"
$ws.Cells.Item(7, 4).Value = "126.8200"
$ws.Cells.Item(7, 5).Value = "128.1400"
$ws.Cells.Item(7, 6).Value = "1274706"

$ws.Cells.Item(8, 1).Value = "2022-09-19"
$ws.Cells.Item(8, 2).Value = "128.7900"
$ws.Cells.Item(8, 3).Value = "130.7900"
$ws.Cells.Item(8, 4).Value = "128.3600"
$ws.Cells.Item(8, 5).Value = "130.7200"
$ws.Cells.Item(8, 6).Value = "1059566"

$ws.Cells.Item(9, 1).Value = "2022-09-16"
$ws.Cells.Item(9, 2).Value = "132.1200"
$ws.Cells.Item(9, 3).Value = "132.1300"
$ws.Cells.Item(9, 4).Value = "128.2800"
$ws.Cells.Item(9, 5).Value = "129.8900"
$ws.Cells.Item(9, 6).Value = "2300604"

$ws.Cells.Item(10, 1).Value = "2022-09-15"
$ws.Cells.Item(10, 2).Value = "133.8100"
$ws.Cells.Item(10, 3).Value = "135.5100"
$ws.Cells.Item(10, 4).Value = "132.6550"
$ws.Cells.Item(10, 5).Value = "133.1600"
$ws.Cells.Item(10, 6).Value = "1446479"

$ws.Cells.Item(11, 1).Value = "2022-09-14"
$ws.Cells.Item(11, 2).Value = "133.8000"
$ws.Cells.Item(11, 3).Value = "135.3600"
$ws.Cells.Item(11, 4).Value = "132.0950"
$ws.Cells.Item(11, 5).Value = "133.2500"
$ws.Cells.Item(11, 6).Value = "1577007"

$ws.Cells.Item(12, 1).Value = "2022-09-13"
$ws.Cells.Item(12, 2).Value = "136.3200"
$ws.Cells.Item(12, 3).Value = "137.4100"
$ws.Cells.Item(12, 4).Value = "133.4400"
$ws.Cells.Item(12, 5).Value = "133.5400"
$ws.Cells.Item(12, 6).Value = "1586967"

# 4) Restore the exact cell styling used by the rest of the data table (date
#    column centered/bold/bordered, value columns unstyled) by copying the
#    formats from row 13 (the first untouched, correctly-styled data row)
#    down across the new rows. This also clears the temporary "@" text
#    format we applied above, matching the original sheet's formatting.
$ws.Range("A13:F13").Copy() | Out-Null
$ws.Range("A2:F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
